$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price/volume refresh, and two row re-orderings)
# Each target value is written with a leading apostrophe to force Excel to
# store it as literal text (matching the source inlineStr cells) instead of
# auto-converting numeric-looking strings to numbers; the style is then reset
# to Normal so no stray quote-prefix formatting / style index is introduced.

$ws.Range("D2").Value2 = "'51.829.10"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value2 = "'  -0.15%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value2 = "'2.790.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value2 = "'  +0.14%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value2 = "'  +0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value2 = "'353.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "'  -1.50%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value2 = "'108.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "'  -0.36%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value2 = "'0.549"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "'  -2.78%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value2 = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "'  +0.09%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value2 = "'0.596"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "'  -0.04%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value2 = "'39.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "'  -0.61%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value2 = "'  +3.03%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value2 = "'20.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value2 = "'  +3.24%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value2 = "'0.0838"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "'  -1.82%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value2 = "'  +0.81%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value2 = "'3.226.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "'  +0.24%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value2 = "'2.760.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "'  -1.70%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value2 = "'0.928"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "'  -0.78%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value2 = "'51.779.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "'  -0.08%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value2 = "'7.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "'  +3.53%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value2 = "'3.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "'  -0.12%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value2 = "'13.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value2 = "'  +0.96%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value2 = "'0.0₃0966"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "'  -1.46%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value2 = "'70.06"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "'  -0.27%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value2 = "'267.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value2 = "'  -2.62%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value2 = "'2.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "'  -0.34%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("B26").Value2 = "'EthereumClassic"
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").Value2 = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").Value2 = "'26.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "'  -1.96%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("B27").Value2 = "'Dai"
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value2 = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value2 = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "'  +0.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value2 = "'0.163"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "'  +12.25%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value2 = "'10.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "'  +0.88%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("B30").Value2 = "'InjectiveProtocol"
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value2 = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value2 = "'36.62"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "'  +7.65%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("B31").Value2 = "'Filecoin"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value2 = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value2 = "'6.25"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "'  +9.51%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value2 = "'Toncoin"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value2 = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value2 = "'2.12"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "'  -7.62%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value2 = "'52.06"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "'  +1.04%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value2 = "'0.0453"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "'  -2.78%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value2 = "'5.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "'  +5.65%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value2 = "'0.0831"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "'  -1.53%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value2 = "'0.997"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value2 = "'  -0.33%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value2 = "'18.58"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "'  +2.64%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value2 = "'3.16"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "'  -2.50%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value2 = "'  -1.46%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value2 = "'2.58"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "'  +0.71%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value2 = "'  -0.62%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value2 = "'120.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value2 = "'22.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "'  +0.31%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value2 = "'  -2.37%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value2 = "'2.122.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value2 = "'  +2.39%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value2 = "'3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "'  +1.30%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value2 = "'  +6.34%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value2 = "'5.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "'  -4.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value2 = "'0.910"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "'  -2.64%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value2 = "'1.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "'  +9.66%  "
$ws.Range("E51").Style = "Normal"
